# Updated cryptos list on Sat Jul 27 06:39:55 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the crypto
# tracker sheet with the latest scraped snapshot, and reorders the three
# rows whose rank changed (Bittensor jumped above InjectiveProtocol and
# Hedera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while always keeping it stored as TEXT
# (matches source data which is inlineStr even for numeric-looking prices).
# A leading apostrophe is Excel's standard "treat as text" quote-prefix;
# harmless for values that are already unambiguous text (e.g. "67.952.57",
# which contains two dots and could never be parsed as a number anyway).
function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
}

# ---------------------------------------------------------------------
# Price (D) column updates
# ---------------------------------------------------------------------
Set-TextValue "D2"  "67.952.57"
Set-TextValue "D3"  "3.259.07"
Set-TextValue "D5"  "585.74"
Set-TextValue "D6"  "183.96"
Set-TextValue "D12" "3.825.45"
Set-TextValue "D15" "67.980.56"
Set-TextValue "D17" "3.255.49"
Set-TextValue "D19" "13.60"
Set-TextValue "D20" "381.03"
Set-TextValue "D22" "1.00"
Set-TextValue "D23" "71.29"
Set-TextValue "D26" "9.84"
Set-TextValue "D30" "5.68"
Set-TextValue "D31" "7.28"
Set-TextValue "D32" "22.88"
Set-TextValue "D33" "0.999"
Set-TextValue "D36" "163.00"
Set-TextValue "D39" "6.79"
Set-TextValue "D40" "26.61"
Set-TextValue "D43" "41.42"
Set-TextValue "D47" "2.644.68"
Set-TextValue "D48" "0.0285"
Set-TextValue "D49" "31.98"

# ---------------------------------------------------------------------
# Volume(1h) (E) column updates
# ---------------------------------------------------------------------
$ws.Range("E2").Value  = "  +1.40%  "
$ws.Range("E3").Value  = "  -0.22%  "
$ws.Range("E5").Value  = "  +1.07%  "
$ws.Range("E6").Value  = "  +3.80%  "
$ws.Range("E7").Value  = "  -0.04%  "
$ws.Range("E8").Value  = "  -1.05%  "
$ws.Range("E9").Value  = "  +3.57%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  +6.88%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E39").Value = "  +5.49%  "
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("E41").Value = "  +5.31%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("E49").Value = "  +5.02%  "
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("E51").Value = "  +1.12%  "

# ---------------------------------------------------------------------
# Rows 44-46 rotated: Bittensor moved up to rank 44 (with refreshed
# price/volume), InjectiveProtocol dropped to 45, Hedera dropped to 46.
# ---------------------------------------------------------------------
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D44" "346.36"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "25.40"
$ws.Range("E45").Value = "  +2.03%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0687"
$ws.Range("E46").Value = "  +1.42%  "
